$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the style of an existing header cell (e.g. F1) onto the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update slightly-changed values in row 2
$ws.Range("B2").Value = 0.02950307763024243
$ws.Range("D2").Value = 0.1292627146720763

# New data cells
$ws.Range("G2").Value = 0.1260932844166139
$ws.Range("H2").Value = 0.991
